$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cellRef, $text)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-CellText "D2" "64.792.75"
Set-CellText "E2" "  +4.01%  "
Set-CellText "D3" "3.089.53"
Set-CellText "E3" "  +2.10%  "
Set-CellText "E4" "  +0.09%  "
Set-CellText "D5" "559.69"
Set-CellText "E5" "  +3.28%  "
Set-CellText "D6" "143.32"
Set-CellText "E6" "  +7.29%  "
Set-CellText "E7" "  +0.02%  "
Set-CellText "D8" "3.082.99"
Set-CellText "E8" "  +2.03%  "
Set-CellText "D9" "0.498"
Set-CellText "E9" "  +1.19%  "
Set-CellText "D10" "6.46"
Set-CellText "E10" "  +4.99%  "
Set-CellText "D11" "0.152"
Set-CellText "E11" "  +3.14%  "
Set-CellText "D12" "0.469"
Set-CellText "E12" "  +5.13%  "
Set-CellText "D13" "0.0000229"
Set-CellText "E13" "  +3.12%  "
Set-CellText "D14" "35.15"
Set-CellText "E14" "  +2.59%  "
Set-CellText "D15" "3.610.66"
Set-CellText "E15" "  +2.85%  "
Set-CellText "D16" "64.867.16"
Set-CellText "E16" "  +4.21%  "
Set-CellText "D17" "3.100.80"
Set-CellText "E17" "  +2.56%  "
Set-CellText "E18" "  +1.03%  "
Set-CellText "D19" "6.74"
Set-CellText "E19" "  +1.76%  "
Set-CellText "D20" "478.09"
Set-CellText "E20" "  -0.32%  "
Set-CellText "D21" "13.73"
Set-CellText "E21" "  +3.65%  "
Set-CellText "D22" "0.676"
Set-CellText "E22" "  +0.69%  "
Set-CellText "D23" "7.54"
Set-CellText "D24" "13.34"
Set-CellText "E24" "  +10.28%  "
Set-CellText "D25" "80.99"
Set-CellText "E25" "  +0.38%  "
Set-CellText "D26" "0.999"
Set-CellText "E26" "  -0.08%  "
Set-CellText "D27" "2.78"
Set-CellText "E27" "  +2.61%  "
Set-CellText "D28" "8.18"
Set-CellText "E28" "  +6.06%  "
Set-CellText "D29" "2.05"
Set-CellText "E29" "  +6.01%  "
Set-CellText "D30" "1.00"
Set-CellText "E30" "  +0.17%  "
Set-CellText "D31" "26.02"
Set-CellText "E31" "  +1.35%  "
Set-CellText "E32" "  +2.58%  "
Set-CellText "D33" "2.47"
Set-CellText "E33" "  +4.34%  "
Set-CellText "D34" "5.61"
Set-CellText "E34" "  -0.55%  "
Set-CellText "D35" "6.14"
Set-CellText "E35" "  +4.68%  "
Set-CellText "D36" "54.88"
Set-CellText "E36" "  -0.20%  "
Set-CellText "D37" "463.64"
Set-CellText "E37" "  +1.88%  "
Set-CellText "D38" "0.0408"
Set-CellText "E38" "  +5.90%  "
Set-CellText "D39" "0.0830"
Set-CellText "E39" "  +3.73%  "
Set-CellText "D40" "2.91"
Set-CellText "E40" "  +18.35%  "
Set-CellText "D41" "2.980.40"
Set-CellText "E41" "  -5.83%  "
Set-CellText "D42" "8.20"
Set-CellText "E42" "  +1.43%  "
Set-CellText "D43" "0.115"
Set-CellText "E43" "  -2.23%  "
Set-CellText "D44" "27.88"
Set-CellText "E44" "  +5.75%  "
Set-CellText "D45" "0.258"
Set-CellText "E45" "  +5.91%  "
Set-CellText "B46" "USDe"
Set-CellText "C46" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-CellText "D46" "1.00"
Set-CellText "E46" "  +0.01%  "
Set-CellText "B47" "Fetch.AI"
Set-CellText "C47" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-CellText "D47" "2.13"
Set-CellText "E47" "  +8.54%  "
Set-CellText "E48" "  +2.80%  "
Set-CellText "D49" "0.0₃0523"
Set-CellText "E49" "  +5.35%  "
Set-CellText "D50" "116.00"
Set-CellText "E50" "  +1.54%  "
Set-CellText "D51" "2.06"
Set-CellText "E51" "  +1.87%  "

$wb.Save()
